# Scheduled market-data refresh for Unicorn_Profits workbook.
# Pulls updated Universalis price snapshots into the currentAveragePrice*
# / LevePrice* / LeveProfit* columns (H:N) for the affected leve rows on each
# job sheet. Values only -- no structural/formatting changes.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 8241.464
$ws.Range("I76").Value = 14046.1
$ws.Range("J76").Value = 5016.6665
$ws.Range("K76").Value = 14046.1
$ws.Range("L76").Value = 5016.6665
$ws.Range("M76").Value = -13731.1
$ws.Range("N76").Value = -5646.6665
$ws.Range("H79").Value = 8241.464
$ws.Range("I79").Value = 14046.1
$ws.Range("J79").Value = 5016.6665
$ws.Range("K79").Value = 14046.1
$ws.Range("L79").Value = 5016.6665
$ws.Range("M79").Value = -12954.1
$ws.Range("N79").Value = -7200.6665
$ws.Range("H113").Value = 3488.3125
$ws.Range("I113").Value = 3115.5557
$ws.Range("J113").Value = 3634.1738
$ws.Range("K113").Value = 3115.5557
$ws.Range("L113").Value = 3634.1738
$ws.Range("M113").Value = 138.4443000000001
$ws.Range("N113").Value = -10142.1738
$ws.Range("H132").Value = 3551.7188
$ws.Range("I132").Value = 2258.913
$ws.Range("K132").Value = 6776.739
$ws.Range("M132").Value = -4246.739

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 384184.38
$ws.Range("I61").Value = 325156.8
$ws.Range("J61").Value = 479752.8
$ws.Range("K61").Value = 325156.8
$ws.Range("L61").Value = 479752.8
$ws.Range("M61").Value = -324944.8
$ws.Range("N61").Value = -480176.8
$ws.Range("H62").Value = 30683
$ws.Range("J62").Value = 30683
$ws.Range("L62").Value = 30683
$ws.Range("N62").Value = -31931
$ws.Range("H63").Value = 3184.6155
$ws.Range("I63").Value = 3040
$ws.Range("J63").Value = 3666.6667
$ws.Range("K63").Value = 3040
$ws.Range("L63").Value = 3666.6667
$ws.Range("M63").Value = -2354
$ws.Range("N63").Value = -5038.6667
$ws.Range("H65").Value = 30683
$ws.Range("J65").Value = 30683
$ws.Range("L65").Value = 92049
$ws.Range("N65").Value = -98289
$ws.Range("H66").Value = 3184.6155
$ws.Range("I66").Value = 3040
$ws.Range("J66").Value = 3666.6667
$ws.Range("K66").Value = 15200
$ws.Range("L66").Value = 18333.3335
$ws.Range("M66").Value = -11768
$ws.Range("N66").Value = -25197.3335
$ws.Range("H132").Value = 21315.908
$ws.Range("I132").Value = 27988.795
$ws.Range("J132").Value = 3966.4
$ws.Range("K132").Value = 83966.38499999999
$ws.Range("L132").Value = 11899.2
$ws.Range("M132").Value = -81436.38499999999
$ws.Range("N132").Value = -16959.2
$ws.Range("H136").Value = 384184.38
$ws.Range("I136").Value = 325156.8
$ws.Range("J136").Value = 479752.8
$ws.Range("K136").Value = 975470.3999999999
$ws.Range("L136").Value = 1439258.4
$ws.Range("M136").Value = -972920.3999999999
$ws.Range("N136").Value = -1444358.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 513.3333
$ws.Range("I128").Value = 513.3333
$ws.Range("K128").Value = 1539.9999
$ws.Range("M128").Value = 950.0001
$ws.Range("H134").Value = 3347.2554
$ws.Range("I134").Value = 3104.5625
$ws.Range("K134").Value = 9313.6875
$ws.Range("M134").Value = -6778.6875
$ws.Range("H139").Value = 37942.5
$ws.Range("J139").Value = 37942.5
$ws.Range("L139").Value = 37942.5
$ws.Range("N139").Value = -48222.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 44459.793
$ws.Range("I99").Value = 74034.28999999999
$ws.Range("J99").Value = 3055.5
$ws.Range("K99").Value = 74034.28999999999
$ws.Range("L99").Value = 3055.5
$ws.Range("M99").Value = -72536.28999999999
$ws.Range("N99").Value = -6051.5
$ws.Range("H126").Value = 44459.793
$ws.Range("I126").Value = 74034.28999999999
$ws.Range("J126").Value = 3055.5
$ws.Range("K126").Value = 222102.87
$ws.Range("L126").Value = 9166.5
$ws.Range("M126").Value = -219632.87
$ws.Range("N126").Value = -14106.5
$ws.Range("H132").Value = 2010.0212
$ws.Range("I132").Value = 950.90625
$ws.Range("J132").Value = 4269.467
$ws.Range("K132").Value = 2852.71875
$ws.Range("L132").Value = 12808.401
$ws.Range("M132").Value = -322.71875
$ws.Range("N132").Value = -17868.401
$ws.Range("H134").Value = 1426.0785
$ws.Range("I134").Value = 998.05
$ws.Range("J134").Value = 2982.5454
$ws.Range("K134").Value = 2994.15
$ws.Range("L134").Value = 8947.636200000001
$ws.Range("M134").Value = -459.1499999999996
$ws.Range("N134").Value = -14017.6362

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 667.3333
$ws.Range("I113").Value = 657
$ws.Range("J113").Value = 750
$ws.Range("K113").Value = 1971
$ws.Range("L113").Value = 2250
$ws.Range("M113").Value = 199
$ws.Range("N113").Value = -6590

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3203.9092
$ws.Range("I102").Value = 1513.6818
$ws.Range("K102").Value = 1513.6818
$ws.Range("M102").Value = 108.3181999999999
$ws.Range("H126").Value = 8800
$ws.Range("J126").Value = 7600
$ws.Range("L126").Value = 22800
$ws.Range("N126").Value = -27740
$ws.Range("H132").Value = 2987.1785
$ws.Range("I132").Value = 2773.15
$ws.Range("J132").Value = 3522.25
$ws.Range("K132").Value = 8319.450000000001
$ws.Range("L132").Value = 10566.75
$ws.Range("M132").Value = -5789.450000000001
$ws.Range("N132").Value = -15626.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 13000
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 14500
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 14500
$ws.Range("M62").Value = -9376
$ws.Range("N62").Value = -15748
$ws.Range("H64").Value = 44500
$ws.Range("J64").Value = 44500
$ws.Range("L64").Value = 44500
$ws.Range("N64").Value = -44950
$ws.Range("H65").Value = 13000
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 14500
$ws.Range("K65").Value = 30000
$ws.Range("L65").Value = 43500
$ws.Range("M65").Value = -26880
$ws.Range("N65").Value = -49740
$ws.Range("H67").Value = 44500
$ws.Range("J67").Value = 44500
$ws.Range("L67").Value = 44500
$ws.Range("N67").Value = -46060
$ws.Range("H122").Value = 3233
$ws.Range("I122").Value = 3233
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9699
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7249
$ws.Range("N122").ClearContents()
$ws.Range("H139").Value = 33187.5
$ws.Range("J139").Value = 33187.5
$ws.Range("L139").Value = 33187.5
$ws.Range("N139").Value = -43467.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 42157.855
$ws.Range("J139").Value = 42157.855
$ws.Range("L139").Value = 42157.855
$ws.Range("N139").Value = -52437.855

